$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-14 Sunday", "2024-07-15 Monday"),
    @("430×2=860", "780×8=6240"),
    @("943×2=1886", "521×5=2605"),
    @("647×8=5176", "574×6=3444"),
    @("741×5=3705", "403×5=2015"),
    @("352×5=1760", "401×9=3609"),
    @("617×3=1851", "763×9=6867"),
    @("598×3=1794", "844×9=7596"),
    @("314×2=628", "480×5=2400"),
    @("504×3=1512", "305×9=2745"),
    @("144×2=288", "700×3=2100"),
    @("861×3=2583", "762×8=6096"),
    @("551×3=1653", "415×7=2905"),
    @("696×5=3480", "358×4=1432"),
    @("942×5=4710", "671×3=2013"),
    @("388×8=3104", "572×7=4004"),
    @("945×6=5670", "229×3=687"),
    @("159×3=477", "462×9=4158"),
    @("760×6=4560", "114×6=684"),
    @("242×4=968", "244×8=1952"),
    @("769×3=2307", "769×8=6152"),
    @("191×8=1528", "621×9=5589"),
    @("811×8=6488", "436×8=3488"),
    @("777×3=2331", "895×9=8055"),
    @("705×2=1410", "982×5=4910"),
    @("534×7=3738", "599×6=3594")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
